{"js": "// Update the two-digit \u00f7 one-digit practice table: replace each\n// \"old\u00f7d=\" expression with its new value, in document order. A couple of\n// source expressions (e.g. \"22\u00f73=\") repeat more than once in the sheet\n// but map to different replacements at each occurrence, so replacements\n// are resolved by matching the Nth search hit to the Nth occurrence of\n// that old value in the mapping (search results come back in document\n// order, same order the occurrences appear top-to-bottom/left-to-right).\nconst replacements = {\n  \"68\u00f72=\": [\"67\u00f72=\"],\n  \"12\u00f79=\": [\"68\u00f77=\"],\n  \"38\u00f77=\": [\"88\u00f76=\"],\n  \"14\u00f73=\": [\"48\u00f79=\"],\n  \"45\u00f76=\": [\"76\u00f78=\"],\n  \"53\u00f76=\": [\"76\u00f76=\"],\n  \"18\u00f72=\": [\"93\u00f76=\"],\n  \"49\u00f74=\": [\"31\u00f72=\"],\n  \"56\u00f78=\": [\"63\u00f75=\"],\n  \"88\u00f75=\": [\"47\u00f78=\"],\n  \"22\u00f73=\": [\"67\u00f77=\", \"69\u00f76=\"],\n  \"67\u00f76=\": [\"78\u00f79=\"],\n  \"34\u00f73=\": [\"79\u00f76=\"],\n  \"79\u00f79=\": [\"84\u00f78=\"],\n  \"29\u00f73=\": [\"66\u00f74=\"],\n  \"78\u00f75=\": [\"77\u00f77=\"],\n  \"23\u00f75=\": [\"31\u00f76=\"],\n  \"87\u00f72=\": [\"30\u00f79=\"],\n  \"49\u00f78=\": [\"19\u00f74=\"],\n  \"77\u00f72=\": [\"98\u00f75=\"],\n  \"67\u00f79=\": [\"99\u00f78=\"],\n  \"84\u00f73=\": [\"85\u00f78=\"],\n  \"58\u00f79=\": [\"63\u00f72=\"],\n  \"52\u00f78=\": [\"98\u00f77=\"]\n};\n\nconst body = context.document.body;\nconst searchResults = {};\nfor (const oldText of Object.keys(replacements)) {\n  searchResults[oldText] = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  searchResults[oldText].load(\"items\");\n}\nawait context.sync();\n\nfor (const [oldText, newTexts] of Object.entries(replacements)) {\n  const items = searchResults[oldText].items;\n  for (let i = 0; i < newTexts.length; i++) {\n    items[i].insertText(newTexts[i], Word.InsertLocation.replace);\n  }\n}\nawait context.sync();\n", "ps1": "# Update the two-digit \u00f7 one-digit practice table: replace each\n# \"old\u00f7d=\" expression with its new value, in document order. A couple of\n# source expressions (e.g. \"22\u00f73=\") repeat more than once in the sheet but\n# map to different replacements at each occurrence, so we walk a single\n# Range across the whole body and do one wdReplaceOne Find per pair, in\n# document order; each Execute() leaves the range positioned just after\n# the replacement, so the next search naturally picks up the next\n# occurrence instead of re-matching an earlier one.\n$d = $word.ActiveDocument\n\n$pairs = @(\n    @(\"68\u00f72=\", \"67\u00f72=\"),\n    @(\"12\u00f79=\", \"68\u00f77=\"),\n    @(\"38\u00f77=\", \"88\u00f76=\"),\n    @(\"14\u00f73=\", \"48\u00f79=\"),\n    @(\"45\u00f76=\", \"76\u00f78=\"),\n    @(\"53\u00f76=\", \"76\u00f76=\"),\n    @(\"18\u00f72=\", \"93\u00f76=\"),\n    @(\"49\u00f74=\", \"31\u00f72=\"),\n    @(\"56\u00f78=\", \"63\u00f75=\"),\n    @(\"88\u00f75=\", \"47\u00f78=\"),\n    @(\"22\u00f73=\", \"67\u00f77=\"),\n    @(\"67\u00f76=\", \"78\u00f79=\"),\n    @(\"34\u00f73=\", \"79\u00f76=\"),\n    @(\"79\u00f79=\", \"84\u00f78=\"),\n    @(\"29\u00f73=\", \"66\u00f74=\"),\n    @(\"78\u00f75=\", \"77\u00f77=\"),\n    @(\"23\u00f75=\", \"31\u00f76=\"),\n    @(\"87\u00f72=\", \"30\u00f79=\"),\n    @(\"49\u00f78=\", \"19\u00f74=\"),\n    @(\"77\u00f72=\", \"98\u00f75=\"),\n    @(\"22\u00f73=\", \"69\u00f76=\"),\n    @(\"67\u00f79=\", \"99\u00f78=\"),\n    @(\"84\u00f73=\", \"85\u00f78=\"),\n    @(\"58\u00f79=\", \"63\u00f72=\"),\n    @(\"52\u00f78=\", \"98\u00f77=\")\n)\n\n$range = $d.Content\n\nforeach ($pair in $pairs) {\n    $range.Find.ClearFormatting()\n    $range.Find.Replacement.ClearFormatting()\n    $range.Find.Text = $pair[0]\n    $range.Find.Replacement.Text = $pair[1]\n    $range.Find.Execute([ref]$pair[0], [ref]$true, [ref]$false, [ref]$false, [ref]$false, [ref]$false, [ref]$true, [ref]1, [ref]$false, [ref]$pair[1], [ref]1)\n}\n"}
